# "Finished the Great Refactoring #1 of Integration Tests"
#
# On the "Generic Backlog" sheet, under the "Professional" section:
#   - Add a new in-progress item: "Replace cooling fan and heat sink"
# And under the "Other Stuff" section:
#   - Remove the completed item "RETURN GI"
#   - Rename "Goto Fleetfeet" to "Gi + Gloves"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generic Backlog")

# --- Add "Replace cooling fan and heat sink" as a new IN PROGRESS row ---
# It belongs right under "Pleiades (...)" (row 7), above "Buy Resharper..." (row 8),
# so insert a fresh row at position 8 and push the rest of the section down.
$ws.Rows(8).Insert()

# Match the formatting of the row above (the other "IN PROGRESS" item).
$ws.Range("A7:B7").Copy()
$ws.Range("A8:B8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A8").Value = "Replace cooling fan and heat sink"
$ws.Range("B8").Value = "IN PROGRESS"

# --- Remove "RETURN GI" entirely ---
# After the insert above, it now lives at row 19 (was row 18).
$ws.Rows(19).Delete()

# --- Rename "Goto Fleetfeet" -> "Gi + Gloves" ---
# After the insert (+1) and the delete above (-1), it's still at row 19.
$ws.Range("A19").Value = "Gi + Gloves"

# Reflect the new active cell/selection on the sheet.
$ws.Range("A19").Select()
